$p = $ppt.ActivePresentation
Write-Output "HasHandoutMaster: $($p.HasHandoutMaster)"
try {
  $hm = $p.HandoutMaster
  Write-Output "hm: $hm"
  $cs = $hm.ColorScheme
  Write-Output "cs: $cs"
} catch { Write-Output "ERR: $_" }
